# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (want-to-go count) figures, and marks one
# event as sold out ("不可售") in its "最低票价" (lowest price) column,
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6994
$wsExhibit.Range("F4").Value = 459
$wsExhibit.Range("G4").Value = "不可售"
$wsExhibit.Range("F7").Value = 132
$wsExhibit.Range("F11").Value = 50
$wsExhibit.Range("F16").Value = 43
$wsExhibit.Range("F17").Value = 3602
$wsExhibit.Range("F21").Value = 19
$wsExhibit.Range("F23").Value = 2227
$wsExhibit.Range("F24").Value = 13
$wsExhibit.Range("F25").Value = 239
$wsExhibit.Range("F31").Value = 155
$wsExhibit.Range("F32").Value = 193
$wsExhibit.Range("F33").Value = 78

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6994
$wsAll.Range("F4").Value = 459
$wsAll.Range("G4").Value = "不可售"
$wsAll.Range("F8").Value = 132
$wsAll.Range("F12").Value = 50
$wsAll.Range("F17").Value = 43
$wsAll.Range("F18").Value = 3602
$wsAll.Range("F22").Value = 19
$wsAll.Range("F24").Value = 2227
$wsAll.Range("F25").Value = 13
$wsAll.Range("F26").Value = 239
$wsAll.Range("F32").Value = 155
$wsAll.Range("F33").Value = 194
$wsAll.Range("F34").Value = 78
